$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 500
$ws.Range("B3").Value = 300
$ws.Range("B4").Value = 200
$ws.Range("B5").Value = 145
$ws.Range("B6").Value = 95
$ws.Range("B7").Value = 300
$ws.Range("B8").Value = 110
